$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 1
$ws.Range("A4").Value = 1
$ws.Range("A7").Value = 2
$ws.Range("A11").Value = 1
$ws.Range("A13").Value = 2
$ws.Range("A15").Value = 1
$ws.Range("A16").Value = 2
$ws.Range("A18").Value = 2
$ws.Range("A21:A23").Value = 1
$ws.Range("A27:A28").Value = 2
$ws.Range("A34:A35").Value = 1
$ws.Range("A38").Value = 2
$ws.Range("A41").Value = 2
$ws.Range("A45:A46").Value = 1
$ws.Range("A48").Value = 1
$ws.Range("A52:A54").Value = 2
$ws.Range("A55").Value = 1
$ws.Range("A58").Value = 1
$ws.Range("A65").Value = 2
$ws.Range("A67").Value = 2
$ws.Range("A71").Value = 1
$ws.Range("A74").Value = 1
$ws.Range("A77").Value = 1
$ws.Range("A80").Value = 1
$ws.Range("A81").Value = 2
$ws.Range("A83").Value = 1
$ws.Range("A85:A87").Value = 2
$ws.Range("A88").Value = 1
$ws.Range("A91").Value = 2
$ws.Range("A93").Value = 1
$ws.Range("A98").Value = 2
$ws.Range("A101:A102").Value = 2
$ws.Range("A104:A105").Value = 1
$ws.Range("A108").Value = 1
$ws.Range("A109").Value = 2
$ws.Range("A111").Value = 1
$ws.Range("A113").Value = 1
$ws.Range("A114").Value = 2
$ws.Range("A116").Value = 1
$ws.Range("A117").Value = 2
$ws.Range("A124").Value = 1
$ws.Range("A126:A127").Value = 1
$ws.Range("A129:A130").Value = 2
$ws.Range("A131").Value = 1
$ws.Range("A133").Value = 1
$ws.Range("A138").Value = 2
$ws.Range("A144").Value = 1
$ws.Range("A146").Value = 2
$ws.Range("A147").Value = 1
$ws.Range("A151:A152").Value = 2
$ws.Range("A154").Value = 1
$ws.Range("A155").Value = 2
$ws.Range("A156:A157").Value = 1
$ws.Range("A159").Value = 1
$ws.Range("A161").Value = 2
$ws.Range("A164").Value = 2
$ws.Range("A166").Value = 2
$ws.Range("A168").Value = 2
$ws.Range("A171").Value = 2
$ws.Range("A172").Value = 1
$ws.Range("A176:A177").Value = 2
$ws.Range("A179").Value = 2
$ws.Range("A181:A182").Value = 1
$ws.Range("A183").Value = 2
$ws.Range("A186").Value = 2
$ws.Range("A187").Value = 1
$ws.Range("A189").Value = 1
$ws.Range("A191").Value = 2
$ws.Range("A196").Value = 2
$ws.Range("A198:A200").Value = 2
